# Update Metrics sheet values (B2:B13) and the active cell selections,
# mirroring a manual edit of the "Metrics" sheet followed by a look at
# the "today" sheet (which holds formulas referencing Metrics!...).
$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Activate()
$metrics.Range("B2").Value  = 419387.58
$metrics.Range("B3").Value  = 318721.65000000008
$metrics.Range("B4").Value  = 111384.94
$metrics.Range("B5").Value  = 17078
$metrics.Range("B6").Value  = 419387.58
$metrics.Range("B7").Value  = 318721.65000000008
$metrics.Range("B8").Value  = 111384.94
$metrics.Range("B9").Value  = 17078
$metrics.Range("B10").Value = 34520639.299999997
$metrics.Range("B11").Value = 32364714.440000001
$metrics.Range("B12").Value = 12057198.799999999
$metrics.Range("B13").Value = 1334985

# Selection in Metrics moved to E14
$metrics.Range("E14").Select()

# The "today" sheet recalculates its formulas that reference Metrics!...
# cells automatically once the precedent values above change; TODAY()-1
# reflects the runtime clock. Just update the selection, matching the
# diff, and leave it as the active (tab-selected) sheet, as it was
# originally.
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("D6").Select()
